$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''37.623.71'
$ws.Range('E2').Value = '  +1.68%  '

# Row 3
$ws.Range('D3').Value = '''2.080.61'
$ws.Range('E3').Value = '  +4.35%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = '''236.71'
$ws.Range('E5').Value = '  -2.96%  '

# Row 6
$ws.Range('D6').Value = '''0.616'
$ws.Range('E6').Value = '  +1.73%  '

# Row 7
$ws.Range('D7').Value = '''58.40'
$ws.Range('E7').Value = '  +6.79%  '

# Row 8
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('D9').Value = '''0.386'
$ws.Range('E9').Value = '  +3.05%  '

# Row 10
$ws.Range('D10').Value = '''58.28'
$ws.Range('E10').Value = '  +1.85%  '

# Row 11
$ws.Range('E11').Value = '  +1.22%  '

# Row 12
$ws.Range('D12').Value = '''0.102'
$ws.Range('E12').Value = '  +3.79%  '

# Row 13
$ws.Range('D13').Value = '''2.390.96'
$ws.Range('E13').Value = '  +4.40%  '

# Row 14
$ws.Range('D14').Value = '''14.54'
$ws.Range('E14').Value = '  +2.66%  '

# Row 15
$ws.Range('D15').Value = '''21.04'
$ws.Range('E15').Value = '  +0.54%  '

# Row 16
$ws.Range('D16').Value = '''0.781'
$ws.Range('E16').Value = '  +2.96%  '

# Row 17
$ws.Range('E17').Value = '  +4.14%  '

# Row 18
$ws.Range('D18').Value = '''2.090.75'
$ws.Range('E18').Value = '  +4.21%  '

# Row 19
$ws.Range('D19').Value = '''37.787.10'
$ws.Range('E19').Value = '  +2.32%  '

# Row 20
$ws.Range('D20').Value = '''6.12'
$ws.Range('E20').Value = '  +20.23%  '

# Row 21
$ws.Range('D21').Value = '''68.70'
$ws.Range('E21').Value = '  +0.15%  '

# Row 22
$ws.Range('E22').Value = '  +0.32%  '

# Row 23
$ws.Range('D23').Value = '''224.76'
$ws.Range('E23').Value = '  -1.69%  '

# Row 24
$ws.Range('E24').Value = '  -0.01%  '

# Row 25
$ws.Range('D25').Value = '''2.42'
$ws.Range('E25').Value = '  +2.95%  '

# Row 26
$ws.Range('E26').Value = '  +0.90%  '

# Row 27
$ws.Range('D27').Value = '''163.25'
$ws.Range('E27').Value = '  +0.46%  '

# Row 28
$ws.Range('D28').Value = '''8.88'
$ws.Range('E28').Value = '  +2.17%  '

# Row 29
$ws.Range('D29').Value = '''0.132'
$ws.Range('E29').Value = '  +4.17%  '

# Row 30
$ws.Range('D30').Value = '''19.38'
$ws.Range('E30').Value = '  +0.85%  '

# Row 31
$ws.Range('E31').Value = '  +7.09%  '

# Row 32
$ws.Range('D32').Value = '''0.118'
$ws.Range('E32').Value = '  +0.61%  '

# Row 33
$ws.Range('D33').Value = '''2.62'
$ws.Range('E33').Value = '  +12.79%  '

# Row 34
$ws.Range('D34').Value = '''4.49'
$ws.Range('E34').Value = '  +1.07%  '

# Row 35
$ws.Range('D35').Value = '''0.0630'
$ws.Range('E35').Value = '  +2.81%  '

# Row 36
$ws.Range('D36').Value = '''4.47'
$ws.Range('E36').Value = '  +5.33%  '

# Row 37
$ws.Range('E37').Value = '  +0.06%  '

# Row 38
$ws.Range('D38').Value = '''5.96'
$ws.Range('E38').Value = '  +12.13%  '

# Row 39
$ws.Range('D39').Value = '''3.34'
$ws.Range('E39').Value = '  +0.23%  '

# Row 40
$ws.Range('E40').Value = '  -0.19%  '

# Row 42
$ws.Range('D42').Value = '''0.0967'
$ws.Range('E42').Value = '  +8.76%  '

# Row 43
$ws.Range('D43').Value = '''1.480.29'
$ws.Range('E43').Value = '  +3.26%  '

# Row 44
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').Value = '''4.32'
$ws.Range('E44').Value = '  +25.71%  '

# Row 45
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '''95.53'
$ws.Range('E45').Value = '  +8.17%  '

# Row 46
$ws.Range('E46').Value = '  +3.64%  '

# Row 47
$ws.Range('D47').Value = '''16.45'
$ws.Range('E47').Value = '  +7.76%  '

# Row 48
$ws.Range('E48').Value = '  +0.54%  '

# Row 49
$ws.Range('D49').Value = '''7.36'
$ws.Range('E49').Value = '  +8.33%  '

# Row 50
$ws.Range('E50').Value = '  +2.18%  '

# Row 51
$ws.Range('D51').Value = '''2.93'
$ws.Range('E51').Value = '  +1.74%  '
